# The document contains several "<id>...</id>" markers that were each
# split across multiple runs (one run per differently-formatted chunk of
# text, e.g. "<id>" / "p106v_" / "1" / "</id>"). This edit collapses two
# of those markers ("p106v_1" and "p107r_1") down into a single run each,
# carrying the text "<id>p106v_1</id>" / "<id>p107r_1</id>" and using the
# formatting of the first (opening-tag) run - matching how Word's
# Find/Replace merges matched text into the first run when a whole match
# is replaced.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p106v_1</id>", $false, $false, $false, $false, $false,
    $true, 1, $false, "<id>p106v_1</id>", 2) | Out-Null

$d.Content.Find.Execute(
    "<id>p107r_1</id>", $false, $false, $false, $false, $false,
    $true, 1, $false, "<id>p107r_1</id>", 2) | Out-Null
